$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '63.535.62'
$ws.Range('E2').Value = '  +3.14%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.070.18'
$ws.Range('E3').Value = '  +2.25%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '551.38'
$ws.Range('E5').Value = '  +2.70%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '142.66'
$ws.Range('E6').Value = '  +6.42%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '3.063.56'
$ws.Range('E8').Value = '  +2.22%  '
$ws.Range('E9').Value = '  +1.35%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '6.54'
$ws.Range('E10').Value = '  +6.71%  '
$ws.Range('E11').Value = '  +2.65%  '
$ws.Range('E12').Value = '  +2.53%  '
$ws.Range('E13').Value = '  +2.72%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '34.96'
$ws.Range('E14').Value = '  +3.38%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.564.80'
$ws.Range('E15').Value = '  +2.24%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '63.489.95'
$ws.Range('E16').Value = '  +3.06%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '3.073.38'
$ws.Range('E17').Value = '  +2.40%  '
$ws.Range('E18').Value = '  -0.95%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '6.79'
$ws.Range('E19').Value = '  +2.88%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '485.83'
$ws.Range('E20').Value = '  +3.94%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '13.88'
$ws.Range('E21').Value = '  +5.10%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.678'
$ws.Range('E22').Value = '  +0.75%  '
$ws.Range('E23').Value = '  +5.59%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '81.13'
$ws.Range('E24').Value = '  +0.80%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '12.79'
$ws.Range('E25').Value = '  +7.13%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range('E27').Value = '  +3.91%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '7.89'
$ws.Range('E28').Value = '  +2.29%  '
$ws.Range('E29').Value = '  +7.65%  '
$ws.Range('E30').Value = '  -0.10%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '26.25'
$ws.Range('E31').Value = '  +2.57%  '
$ws.Range('E32').Value = '  +1.40%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '2.46'
$ws.Range('E33').Value = '  +7.75%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.71'
$ws.Range('E34').Value = '  +4.50%  '
$ws.Range('E35').Value = '  +0.59%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '6.02'
$ws.Range('E36').Value = '  +2.04%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '467.69'
$ws.Range('E37').Value = '  +3.40%  '
$ws.Range('E38').Value = '  +5.14%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.0399'
$ws.Range('E39').Value = '  +4.01%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '3.041.25'
$ws.Range('E40').Value = '  -3.98%  '
$ws.Range('E41').Value = '  -0.73%  '
$ws.Range('E42').Value = '  +2.13%  '
$ws.Range('E43').Value = '  +5.50%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '27.84'
$ws.Range('E44').Value = '  +4.59%  '
$ws.Range('E45').Value = '  +5.37%  '
$ws.Range('E47').Value = '  +2.85%  '
$ws.Range('E48').Value = '  +2.62%  '
$ws.Range('B49').Value = 'PEPE'
$ws.Range('C49').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0₃0512'
$ws.Range('E49').Value = '  +3.45%  '
$ws.Range('B50').Value = 'Monero'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '116.66'
$ws.Range('E50').Value = '  -1.60%  '
$ws.Range('E51').Value = '  +4.60%  '
